# Update practice block fish color assignments (Sheet1, columns C:F, rows 5-16).
# Shared-string pool used in this workbook for these cells:
#   "stimuli/bead_y.PNG"  (yellow)
#   "stimuli/bead_b.PNG"  (blue)
#   "stimuli/bead_g.PNG"  (green)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C5").Value  = "stimuli/bead_b.PNG"

$ws.Range("C6").Value  = "stimuli/bead_g.PNG"
$ws.Range("D6").Value  = "stimuli/bead_b.PNG"

$ws.Range("C7").Value  = "stimuli/bead_b.PNG"
$ws.Range("D7").Value  = "stimuli/bead_g.PNG"
$ws.Range("E7").Value  = "stimuli/bead_b.PNG"

$ws.Range("C8").Value  = "stimuli/bead_g.PNG"
$ws.Range("D8").Value  = "stimuli/bead_b.PNG"
$ws.Range("E8").Value  = "stimuli/bead_g.PNG"
$ws.Range("F8").Value  = "stimuli/bead_b.PNG"

$ws.Range("D9").Value  = "stimuli/bead_g.PNG"
$ws.Range("E9").Value  = "stimuli/bead_b.PNG"
$ws.Range("F9").Value  = "stimuli/bead_g.PNG"

$ws.Range("E10").Value = "stimuli/bead_g.PNG"
$ws.Range("F10").Value = "stimuli/bead_b.PNG"

$ws.Range("C11").Value = "stimuli/bead_y.PNG"
$ws.Range("F11").Value = "stimuli/bead_g.PNG"

$ws.Range("D12").Value = "stimuli/bead_y.PNG"

$ws.Range("E13").Value = "stimuli/bead_y.PNG"

$ws.Range("C14").Value = "stimuli/bead_g.PNG"
$ws.Range("F14").Value = "stimuli/bead_y.PNG"

$ws.Range("C15").Value = "stimuli/bead_y.PNG"
$ws.Range("D15").Value = "stimuli/bead_g.PNG"

$ws.Range("D16").Value = "stimuli/bead_y.PNG"
$ws.Range("E16").Value = "stimuli/bead_g.PNG"

# Reflect the author's final on-screen selection for this sheet.
$ws.Range("C17:F22").Select()
